$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Exceptions")

$ws.Range("A6").Value = "0x0004"
$ws.Range("B6").Value = "DataBase.cs"
$ws.Range("C6").Value = "0x0004, Failed to ReadFromTable,"

$ws.Range("C6").Select()
